$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A10").Value = "Globo"
$ws.Range("B10").Value = "RJ TV 1"
$ws.Range("C10").Value = "Social"
$ws.Range("D10").Value = "2025-04-01T12:36"
$ws.Range("E10").Value = "Positivo"
$ws.Range("F10").Value = "Oportunidades de trabalho. Em Campos, 366 vagas, entre elas para taifeiro e nutricionista offshore, Garçom e auxiliar de serviços gerais.  "
